# Auto-generated Excel COM-interop script
# Updates market-price data cells (columns H-N) across multiple worksheets
# to reflect a refreshed scrape, per the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4740
$ws.Range("I51").Value = 4725
$ws.Range("J51").Value = 4750
$ws.Range("K51").Value = 4725
$ws.Range("L51").Value = 4750
$ws.Range("M51").Value = -4241
$ws.Range("N51").Value = -5718
$ws.Range("H58").Value = 364.23077
$ws.Range("I58").Value = 248.72728
$ws.Range("J58").Value = 999.5
$ws.Range("K58").Value = 746.18184
$ws.Range("L58").Value = 2998.5
$ws.Range("M58").Value = -596.18184
$ws.Range("N58").Value = -3298.5
$ws.Range("H98").Value = 1495.75
$ws.Range("I98").Value = 1512.5454
$ws.Range("J98").Value = 1458.8
$ws.Range("K98").Value = 1512.5454
$ws.Range("L98").Value = 1458.8
$ws.Range("M98").Value = -14.54539999999997
$ws.Range("N98").Value = -4454.8
$ws.Range("H112").Value = 5182.7617
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 5391.9
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 16175.7
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -18391.7
$ws.Range("H122").Value = 1495.75
$ws.Range("I122").Value = 1512.5454
$ws.Range("J122").Value = 1458.8
$ws.Range("K122").Value = 4537.6362
$ws.Range("L122").Value = 4376.4
$ws.Range("M122").Value = -2087.6362
$ws.Range("N122").Value = -9276.4
$ws.Range("H137").Value = 7126.6313
$ws.Range("I137").Value = 2901.75
$ws.Range("J137").Value = 14369.286
$ws.Range("K137").Value = 8705.25
$ws.Range("L137").Value = 43107.858
$ws.Range("M137").Value = -6155.25
$ws.Range("N137").Value = -48207.858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22376.941
$ws.Range("I32").Value = 22717.686
$ws.Range("K32").Value = 22717.686
$ws.Range("M32").Value = -22430.686
$ws.Range("H61").Value = 1391948.6
$ws.Range("I61").Value = 1590300.1
$ws.Range("K61").Value = 1590300.1
$ws.Range("M61").Value = -1590088.1
$ws.Range("H74").Value = 3008.5217
$ws.Range("I74").Value = 1890.6364
$ws.Range("J74").Value = 4033.25
$ws.Range("K74").Value = 1890.6364
$ws.Range("L74").Value = 4033.25
$ws.Range("M74").Value = -1016.6364
$ws.Range("N74").Value = -5781.25
$ws.Range("H77").Value = 3008.5217
$ws.Range("I77").Value = 1890.6364
$ws.Range("J77").Value = 4033.25
$ws.Range("K77").Value = 9453.182000000001
$ws.Range("L77").Value = 20166.25
$ws.Range("M77").Value = -5085.182000000001
$ws.Range("N77").Value = -28902.25
$ws.Range("H97").Value = 1112.591
$ws.Range("I97").Value = 819.8946999999999
$ws.Range("K97").Value = 819.8946999999999
$ws.Range("M97").Value = -323.8946999999999
$ws.Range("H136").Value = 1391948.6
$ws.Range("I136").Value = 1590300.1
$ws.Range("K136").Value = 4770900.300000001
$ws.Range("M136").Value = -4768350.300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1286.1578
$ws.Range("I94").Value = 1202.8235
$ws.Range("K94").Value = 1202.8235
$ws.Range("M94").Value = -751.8235
$ws.Range("H105").Value = 1797.3334
$ws.Range("I105").Value = 1704.6154
$ws.Range("K105").Value = 1704.6154
$ws.Range("M105").Value = 42.38460000000009
$ws.Range("H134").Value = 1776410.9
$ws.Range("I134").Value = 1833747.5
$ws.Range("K134").Value = 5501242.5
$ws.Range("M134").Value = -5498707.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 395.91666
$ws.Range("I7").Value = 503
$ws.Range("J7").Value = 319.42856
$ws.Range("K7").Value = 503
$ws.Range("L7").Value = 319.42856
$ws.Range("M7").Value = -390
$ws.Range("N7").Value = -545.4285600000001
$ws.Range("H22").Value = 1004.7143
$ws.Range("I22").Value = 743
$ws.Range("K22").Value = 743
$ws.Range("M22").Value = -393

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3200.5
$ws.Range("I5").Value = 3701.5
$ws.Range("J5").Value = 2699.5
$ws.Range("K5").Value = 11104.5
$ws.Range("L5").Value = 8098.5
$ws.Range("M5").Value = -10992.5
$ws.Range("N5").Value = -8322.5
$ws.Range("H135").Value = 3200.5
$ws.Range("I135").Value = 3701.5
$ws.Range("J135").Value = 2699.5
$ws.Range("K135").Value = 33313.5
$ws.Range("L135").Value = 24295.5
$ws.Range("M135").Value = -30778.5
$ws.Range("N135").Value = -29365.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 22210.375
$ws.Range("J107").Value = 461.1111
$ws.Range("L107").Value = 461.1111
$ws.Range("N107").Value = -4301.1111
$ws.Range("H113").Value = 5007.486
$ws.Range("I113").Value = 3964.9048
$ws.Range("K113").Value = 3964.9048
$ws.Range("M113").Value = -1794.9048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1121.0625
$ws.Range("I22").Value = 759.7778
$ws.Range("K22").Value = 759.7778
$ws.Range("M22").Value = -464.7778
$ws.Range("H27").Value = 1121.0625
$ws.Range("I27").Value = 759.7778
$ws.Range("K27").Value = 759.7778
$ws.Range("M27").Value = -652.7778
$ws.Range("H46").Value = 1475
$ws.Range("J46").Value = 1299.6666
$ws.Range("L46").Value = 1299.6666
$ws.Range("N46").Value = -1675.6666
$ws.Range("H122").Value = 4614.7856
$ws.Range("I122").Value = 4217.875
$ws.Range("K122").Value = 12653.625
$ws.Range("M122").Value = -10203.625
$ws.Range("H132").Value = 10597.061
$ws.Range("I132").Value = 11025.214
$ws.Range("K132").Value = 33075.642
$ws.Range("M132").Value = -30545.642

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 4360.6
$ws.Range("I3").Value = 1450
$ws.Range("K3").Value = 1450
$ws.Range("M3").Value = -1336
$ws.Range("H11").Value = 10029
$ws.Range("I11").Value = 6748.3335
$ws.Range("J11").Value = 14950
$ws.Range("K11").Value = 6748.3335
$ws.Range("L11").Value = 14950
$ws.Range("M11").Value = -6606.3335
$ws.Range("N11").Value = -15234
$ws.Range("H70").Value = 52951.668
$ws.Range("J70").Value = 52951.668
$ws.Range("L70").Value = 52951.668
$ws.Range("N70").Value = -53581.668
$ws.Range("H73").Value = 52951.668
$ws.Range("J73").Value = 52951.668
$ws.Range("L73").Value = 52951.668
$ws.Range("N73").Value = -55135.668
$ws.Range("H74").Value = 33510.168
$ws.Range("J74").Value = 10212.2
$ws.Range("L74").Value = 10212.2
$ws.Range("N74").Value = -12084.2
$ws.Range("H77").Value = 33510.168
$ws.Range("J77").Value = 10212.2
$ws.Range("L77").Value = 30636.6
$ws.Range("N77").Value = -39996.60000000001
$ws.Range("H107").Value = 1609.5714
$ws.Range("I107").Value = 856.1667
$ws.Range("K107").Value = 2568.5001
$ws.Range("M107").Value = -648.5001000000002
$ws.Range("H132").Value = 2690585.5
$ws.Range("I132").Value = 3032599.8
$ws.Range("J132").Value = 3330
$ws.Range("K132").Value = 9097799.399999999
$ws.Range("L132").Value = 9990
$ws.Range("M132").Value = -9095269.399999999
$ws.Range("N132").Value = -15050
